# "sum 22 week 12 inputs"
# Appends 20 new matchup rows (week 12 results) to the bottom of Sheet1's
# A:D data table (Player_1, Points_1, Player_2, Points_2), then moves the
# saved selection/scroll position to just past the newly entered data -
# mirroring where Excel leaves the cursor after typing in a fresh block
# of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Existing data runs through row 1129 (A1:D1129). New rows are appended
# starting at row 1130.
$startRow = 1130

$newRows = @(
    @(5,1,4,2),
    @(3,1,3,2),
    @(5,0,5,2),
    @(6,3,3,0),
    @(3,2,2,1),
    @(3,2,2,1),
    @(6,0,6,2),
    @(3,3,3,0),
    @(6,2,6,0),
    @(3,1,3,2),
    @(4,2,3,1),
    @(3,3,3,0),
    @(4,3,4,0),
    @(6,2,5,0),
    @(5,2,6,0),
    @(5,2,5,0),
    @(4,0,4,2),
    @(2,3,3,0),
    @(6,3,5,0),
    @(4,0,5,2)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

$lastRow = $startRow + $newRows.Length - 1   # 1149
$nextRow = $lastRow + 1                      # 1150

# Reproduce the saved view state: scrolled so row 1127 is at the top of
# the window, with the active cell sitting one row below the last entry.
$ws.Activate() | Out-Null
$ws.Range("A" + $nextRow).Select() | Out-Null
$excel.Goto($ws.Range("A1127"), $true) | Out-Null
$ws.Range("A" + $nextRow).Select() | Out-Null
